$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("TestData")

# Fix Hobbies value in row 2 (was "Cricket&Movies", now "Cricket")
$ws.Range("I2").Value = "Cricket"

# Fill in Result column (S) with "PASS" for rows 2,3,5,6,8
$ws.Range("S2").Value = "PASS"
$ws.Range("S3").Value = "PASS"
$ws.Range("S5").Value = "PASS"
$ws.Range("S6").Value = "PASS"
$ws.Range("S8").Value = "PASS"

# Update the active selection to J2
$ws.Range("J2").Select()

$wb.Save()
